$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of several K-column cells (the style/cell itself is kept,
# only the inline string value is removed).
$ws.Range("K3").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("K13").ClearContents()
$ws.Range("K14").ClearContents()

# Update the "GELEN SWIFT" maximum amount for AKBANK (column D, row 13)
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
